$d = $word.ActiveDocument

# 1. Update the letter date (unique string in the document).
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the recipient's mailing-address paragraph ("2141 Rancho Mccormick Ct.,
#    Santa Clara CA 95050") into two separate paragraphs: the street line and a
#    new "Santa Clara, CA 95050" line. Only the address block near the top of the
#    letter is affected -- the identical text inside the "PROPERTY ADDRESS:" table
#    further down must stay untouched, so we scope the edit to the paragraph that
#    immediately follows "Anita Blanco".
$addrPara = $d.Paragraphs.Item(7)
$addrPara.Range.InsertParagraphAfter()

$newAddrPara = $d.Paragraphs.Item(8)
$newAddrPara.Range.InsertBefore("Santa Clara, CA 95050")

$d.Paragraphs.Item(7).Range.Find.Execute("2141 Rancho Mccormick Ct., Santa Clara CA 95050", $true, $false, $false, $false, $false,
                   $true, 1, $false, "2141 Rancho Mccormick Ct.", 2)

# 3. Remove the blank "No Spacing" paragraph that used to sit right after the
#    "... Board of Directors" line. Locate it via Find (rather than a fixed
#    paragraph index) so it keeps working even though step 2 shifted every
#    later paragraph index down by one.
$boardRng = $d.Content
$boardRng.Find.Execute("Board of Directors", $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 0)
$boardRng.Expand(4)
$blankRange = $d.Range($boardRng.End, $boardRng.End + 1)
$blankRange.Delete()
